# Apply text corrections (gender/number agreement fixes for
# "Competencias (Portfolio)" -> "Competencia (Portfolio)" singular feminine forms)
# as described by the XML diff, across every cell that uses the affected
# shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# New text values, keyed by the logical change they represent.
$textCadastradas = "SYSTEM exibe a listagem das Competencias (Portfolio) cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$textDestaca      = "SYSTEM destaca a Competencia (Portfolio) selecionada na listagem"
$textExcluirOpcao = "Lider de Pessoas clica na opcao 'Excluir' para excluir a Competencia (Portfolio) selecionada"
$textSemExcluida  = "SYSTEM exibe a listagem das Competencias (Portfolio) sem a Competencia (Portfolio) excluida"
$textComExcluida  = "SYSTEM exibe a listagem das Competencias (Portfolio) com a Competencia (Portfolio) excluida"
$textSoVisualiza  = "SYSTEM exibe a listagem das Competencias (Portfolio) cadastradas apenas para visualizacao com a opcao 'Ajuda'"
$textEditarOpcao  = "Lider de Pessoas clica na opcao 'Editar' para modificar a Competencia (Portfolio) selecionada"

# Cell -> new text map (covers every occurrence of the affected shared strings
# in the worksheet, including the repeated blocks for each test case).
$updates = @{
    "D10"  = $textCadastradas
    "D20"  = $textCadastradas
    "D30"  = $textCadastradas
    "D47"  = $textCadastradas
    "D60"  = $textCadastradas
    "D73"  = $textCadastradas
    "D86"  = $textCadastradas
    "D99"  = $textCadastradas
    "D113" = $textCadastradas
    "D127" = $textCadastradas
    "D142" = $textCadastradas
    "D157" = $textCadastradas
    "D171" = $textCadastradas
    "D185" = $textCadastradas
    "D199" = $textCadastradas
    "D213" = $textCadastradas

    "D11"  = $textDestaca
    "D21"  = $textDestaca
    "D31"  = $textDestaca
    "D128" = $textDestaca
    "D143" = $textDestaca
    "D158" = $textDestaca
    "D172" = $textDestaca
    "D186" = $textDestaca
    "D200" = $textDestaca
    "D214" = $textDestaca

    "B12"  = $textExcluirOpcao
    "B22"  = $textExcluirOpcao
    "B32"  = $textExcluirOpcao

    "D13"  = $textSemExcluida
    "D33"  = $textComExcluida
    "D40"  = $textSoVisualiza

    "B129" = $textEditarOpcao
    "B144" = $textEditarOpcao
    "B159" = $textEditarOpcao
    "B173" = $textEditarOpcao
    "B187" = $textEditarOpcao
    "B201" = $textEditarOpcao
    "B215" = $textEditarOpcao
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
